# "JavaDoc Kommentare erstellt und Testing sheet ausgefüllt"
# Fill in the "Test Design - Use Case based" sheet (Register / Login use
# cases) and append new rows to the "Test Reports" sheet, plus a handful
# of cosmetic layout tweaks (column widths, row heights, zoom, page setup).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Test Design - Use Case based"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Header row
$ws1.Range("A1").Value = "Test Case ID"
$ws1.Range("B1").Value = "Use Case ID"
$ws1.Range("C1").Value = "Flow"
$ws1.Range("D1").Value = "Inputs and additional steps"
$ws1.Range("E1").Value = "Environment"
$ws1.Range("F1").Value = "Expected Result"
$ws1.Range("G1").Value = "Comments"

# Use Case "Register" (test case 1)
$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = "1, Register"
$ws1.Range("C2").Value = "Main Flow"
$ws1.Range("D2").Value = "Max als Benutzername, Meldung mit ok bestätigen"
$ws1.Range("E2").Value = 'Test database (ohne User "Max")'
$ws1.Range("F2").Value = "Erfolgsmeldung mit Benutzername und Passwort, zurück zum Menü"

$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = "1, Register"
$ws1.Range("C3").Value = "Alternative Flow 1"
$ws1.Range("D3").Value = "admin als Benutzername bei erster Abfrage, Meldung mit ok bestätigen, Mustermann bei erneuter Abfrage"
$ws1.Range("E3").Value = 'Test database (mit existierendem User "admin" ohne "Mustermann")'
$ws1.Range("F3").Value = "Fehlermeldung: der Benutzer existiert bereits und Abfrage nach neuem Benutzernamen, zurück zum Menü"
$ws1.Range("G3").Value = "Benutzername existiert bereits"

$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = "1, Register"
$ws1.Range("C4").Value = "Alternative Flow 2"
$ws1.Range("D4").Value = "admin als Benutzername dann abbrechen bei Fehlermeldung"
$ws1.Range("E4").Value = 'Test database (mit existierendem User "admin")'
$ws1.Range("F4").Value = "Fehlermeldung: Benutzer existiert bereits, zurück zum Menü"
$ws1.Range("G4").Value = "Nach Meldung über existierenden User bricht der Benutzer ab"

$ws1.Range("A5").Value = 4
$ws1.Range("B5").Value = "1, Register"
$ws1.Range("C5").Value = "Alternative Flow 3"
# leading "'" -> quote-prefix text cell (matches the workbook's D5 style)
$ws1.Range("D5").Value = "'Abbrechen bei Abfrage auf Benutzernamen und Passwort"
$ws1.Range("E5").Value = "Test database (mit bereits existierenden Usern)"
$ws1.Range("F5").Value = "Benutzer landet wieder im Menü"
$ws1.Range("G5").Value = "Benutzer bricht bei Registriervorgang ab"

# Use Case "Login" (test case 10)
$ws1.Range("A6").Value = 5
$ws1.Range("B6").Value = "10, Login"
$ws1.Range("C6").Value = "Main Flow"
$ws1.Range("D6").Value = "Benutzername und Passwort admin, Meldung mit ok bestätigen"
$ws1.Range("E6").Value = 'Test database (mit existierendem User "admin")'
$ws1.Range("F6").Value = "Erfolgsmeldung über Anmeldung und zurück zum Menü"

$ws1.Range("A7").Value = 6
$ws1.Range("B7").Value = "10, Login"
$ws1.Range("C7").Value = "Alternative Flow 1"
$ws1.Range("D7").Value = "Admin bei erster Abfrage, Meldung mit ok bestätigen, admin bei zweiter Abfrage, Meldung mit ok bestätigen"
$ws1.Range("E7").Value = 'Test database (mit existierendem User "admin")'
$ws1.Range("F7").Value = "Meldung über falsche Daten, erneute Dateneingabe, Erfolgsmeldung, zurück zum Menü"
$ws1.Range("G7").Value = "Benutzer gibt einmal falsche Anmeldedaten ein und die korrekten"

$ws1.Range("A8").Value = 7
$ws1.Range("B8").Value = "10, Login"
$ws1.Range("C8").Value = "Alternative Flow 2"
$ws1.Range("D8").Value = "Abbrechen bei Abfrage auf Benutzernamen und Passwort"
$ws1.Range("E8").Value = "Test database (mit bereits existierenden Usern)"
$ws1.Range("F8").Value = "Benutzer landet wieder im Menü"
$ws1.Range("G8").Value = "Benutzer bricht bei Anmeldevorgang ab"

$ws1.Range("A9").Value = 8
$ws1.Range("B9").Value = "10, Login"
$ws1.Range("C9").Value = "Alternative Flow 3"
$ws1.Range("D9").Value = "Benutzername und Passwort Admin, bei Fehlermeldung abbrechen"
$ws1.Range("E9").Value = 'Test database (ohne User "Admin")'
$ws1.Range("F9").Value = "Fehlermeldung: falsche Anmeldedaten und bei Abbrechen zurück zum Menü"
$ws1.Range("G9").Value = "Benutzer bricht nach falschen Anmeldedaten ab"

# Row heights for the now-much-longer wrapped text (row 2 keeps its
# original 31.5 auto height)
$ws1.Rows.Item(3).RowHeight = 63
$ws1.Rows.Item(4).RowHeight = 51
$ws1.Rows.Item(5).RowHeight = 41.25
$ws1.Rows.Item(6).RowHeight = 51
$ws1.Rows.Item(7).RowHeight = 53.25
$ws1.Rows.Item(8).RowHeight = 42
$ws1.Rows.Item(9).RowHeight = 42

# Column widths (B, C, D, F widened; new G column added)
$ws1.Columns.Item(2).ColumnWidth = 11
$ws1.Columns.Item(3).ColumnWidth = 16.166666666666664
$ws1.Columns.Item(4).ColumnWidth = 91.5
$ws1.Columns.Item(6).ColumnWidth = 33.5
$ws1.Columns.Item(7).ColumnWidth = 31.833333333333336

# View: zoom to 85% and move the selection
$excel.ActiveWindow.Zoom = 85
$ws1.Range("C10").Select() | Out-Null

# Page setup for printing
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Sheet 2: "Test Reports" - log the 8 executed test cases above
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B1").Value = "Time"
$ws2.Range("C1").Value = "Branch, Version"

$reportRows = @(
    @(2, 0.49444444444444446, 1),
    @(3, 0.49513888888888885, 2),
    @(4, 0.49652777777777773, 3),
    @(5, 0.49722222222222223, 4),
    @(6, 0.49791666666666662, 5),
    @(7, 0.49861111111111112, 6),
    @(8, 0.49861111111111112, 7),
    @(9, 0.5, 8)
)

foreach ($entry in $reportRows) {
    $r = $entry[0]
    $ws2.Range("A$r").Value = 45103
    $ws2.Range("B$r").Value = $entry[1]
    $ws2.Range("C$r").Value = "main"
    $ws2.Range("D$r").Value = $entry[2]
    $ws2.Range("E$r").Value = "Kierstein"
    $ws2.Range("F$r").Value = "Success"
    $ws2.Range("G$r").Value = "-"
}

$ws2.Range("C12").Select() | Out-Null
